$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (CPDMA, D2D table)
$ws.Range("C3").Value = 25.945539
$ws.Range("D3").Value = 36.069733
$ws.Range("E3").Value = 40.283211
$ws.Range("F3").Value = 42.803071
$ws.Range("G3").Value = 44.929871
$ws.Range("H3").Value = 48.40069
$ws.Range("I3").Value = 45.222307
$ws.Range("J3").Value = 45.184094
$ws.Range("K3").Value = 45.154422
$ws.Range("L3").Value = 45.158995

# Row 8 (CPDMA, H2D table)
$ws.Range("C8").Value = 6.042957
$ws.Range("D8").Value = 8.599640000000001
$ws.Range("E8").Value = 10.402053
$ws.Range("F8").Value = 11.867319
$ws.Range("G8").Value = 12.655001
$ws.Range("H8").Value = 12.971142
$ws.Range("I8").Value = 13.075969
$ws.Range("J8").Value = 13.465569
$ws.Range("K8").Value = 13.214897
$ws.Range("L8").Value = 13.484385

# Row 13 (CPDMA, D2H table)
$ws.Range("C13").Value = 6.072265
$ws.Range("D13").Value = 8.841113999999999
$ws.Range("E13").Value = 10.579234
$ws.Range("F13").Value = 11.785123
$ws.Range("G13").Value = 12.726844
$ws.Range("H13").Value = 13.291776
$ws.Range("I13").Value = 13.368646
$ws.Range("J13").Value = 13.681402
$ws.Range("K13").Value = 13.620134
$ws.Range("L13").Value = 13.581459
